# "refatorei os codigos de segmentacao, rfma, faturamento e metricas de
#  recorrencia e retencao para bibi"
#
# For the annual recurrence-metrics sheet: a new earliest year (2021) is
# inserted as the first data row, pushing the existing 2022/2023/2024 rows
# down by one row each (their figures are carried over unchanged, only the
# retention_rate of the (now) 2022 row is recomputed against the new prior
# year). A brand-new row for 2025 is appended at the bottom with fresh
# figures, replacing the previous (now stale) 2025 placeholder row.
# Net effect: every data row (2-6) gets rewritten to its final value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-YearRow($Row, $Year, $TotalCustomers, $ReturningCustomers, $NewCustomers, $RetentionRate, $NewRate, $ReturningRate) {
    # ano / ano_obj are stored as text (shared strings) in the source data,
    # not as numbers, even though they look numeric - force text with a
    # leading apostrophe, then clear the auto-applied "quote prefix" number
    # format so the cell keeps the plain default style (no explicit "s").
    $ws.Cells.Item($Row, 1).Value = "'" + $Year
    $ws.Cells.Item($Row, 1).ClearFormats()
    $ws.Cells.Item($Row, 2).Value = "'" + $Year
    $ws.Cells.Item($Row, 2).ClearFormats()

    $ws.Cells.Item($Row, 3).Value = $TotalCustomers
    $ws.Cells.Item($Row, 4).Value = $ReturningCustomers
    $ws.Cells.Item($Row, 5).Value = $NewCustomers
    $ws.Cells.Item($Row, 6).Value = $RetentionRate
    $ws.Cells.Item($Row, 7).Value = $NewRate
    $ws.Cells.Item($Row, 8).Value = $ReturningRate
}

Set-YearRow 2 "2021" 485 66 419 60.55045871559633 86.39175257731959 13.60824742268041
Set-YearRow 3 "2022" 362 167 195 34.43298969072165 53.86740331491713 46.13259668508287
Set-YearRow 4 "2023" 320 186 134 51.38121546961326 41.875 58.12500000000001
Set-YearRow 5 "2024" 457 233 224 72.8125 49.01531728665208 50.98468271334792
Set-YearRow 6 "2025" 340 290 50 63.45733041575492 14.70588235294118 85.29411764705883
